# Update "想去人数" (want-to-go count) figures in column F for the sheets
# that list exhibition rows ("展览" and "全部类型"), matching the refreshed
# gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 648
    5  = 1113
    7  = 11718
    10 = 466
    11 = 385
    14 = 13409
    15 = 13265
    23 = 146
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
